# Applies the indentation tweaks for NORMAM 01/2026 update:
#  - "(Firmas Reconhecidas por semelhanca)" paragraph: LeftIndent 5284 -> 5103 twips
#  - "(Representante da CP/DL/AG)" paragraph: FirstLineIndent 214 -> 72 twips

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text

    if ($text -like "*(Firmas Reconhecidas*") {
        $p.Range.ParagraphFormat.LeftIndent = 5103 / 20
    }
    elseif ($text -like "*(Representante da CP*") {
        $p.Range.ParagraphFormat.FirstLineIndent = 72 / 20
    }
}
